$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Name" column (B), shifting it to C,
# and making room for a "Team_code" column at B.
$ws.Columns.Item(2).Insert() | Out-Null

# Header row
$ws.Range("B1").Value = "Team_code"

# Team_code values per team group
$ws.Range("B2:B14").Value = "admin"
$ws.Range("B15:B21").Value = "support"
$ws.Range("B22:B28").Value = "service"

# name_code header + formula column: "<team_code>_<running count of team_code>"
$ws.Range("D1").Value = "name_code"
for ($r = 2; $r -le 28; $r++) {
    $ws.Range("D$r").Formula = "=B$r&""_""&COUNTIF(`$B`$1:B$r,B$r)"
}

# Give the new D column the same per-row formatting (borders/style) as the
# rest of the table, reusing the existing style indices (copy the format
# from column A of the same row rather than assigning a brand-new style).
for ($r = 1; $r -le 28; $r++) {
    $ws.Range("A$r").Copy() | Out-Null
    $ws.Range("D$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Column widths
$ws.Columns.Item(2).ColumnWidth = 12.29
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

# Selection, as left by the editor
$ws.Range("F5").Select() | Out-Null
